$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 3.722691922189325
$ws.Range("G2").Value = 3.524216379249138
$ws.Range("H2").Value = 3.909279656780932
$ws.Range("I2").Value = 0.654594052840987
$ws.Range("J2").Value = 0.6425483925771267
$ws.Range("K2").Value = 0.6671019205430807
$ws.Range("L2").Value = 0.05083348664753782
$ws.Range("M2").Value = 0.04986780401849426
$ws.Range("N2").Value = 0.05182526186168855

# Row 3
$ws.Range("F3").Value = 0.00322234067449196
$ws.Range("G3").Value = 0.002419159407689076
$ws.Range("H3").Value = 0.004126840667797761
$ws.Range("I3").Value = 0.002986376399492045
$ws.Range("J3").Value = 0.002232978267157508
$ws.Range("K3").Value = 0.003830695812969464
$ws.Range("L3").Value = 0.003364225886023208
$ws.Range("M3").Value = 0.002546671641955229
$ws.Range("N3").Value = 0.004282845019355047

# Row 4
$ws.Range("F4").Value = 3.725914262863816
$ws.Range("G4").Value = 3.526635538656826
$ws.Range("H4").Value = 3.913406497448729
$ws.Range("I4").Value = 0.6575804292404792
$ws.Range("J4").Value = 0.6447813708442841
$ws.Range("K4").Value = 0.6709326163560502
$ws.Range("L4").Value = 0.05419771253356103
$ws.Range("M4").Value = 0.0524144756604495
$ws.Range("N4").Value = 0.0561081068810436
